$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column H (Absent) set to 1 for rows 3-8
$ws.Range("H3:H8").Value = 1

# Row 9: G9 (Invalid) and H9 (Absent) set to 1
$ws.Range("G9").Value = 1
$ws.Range("H9").Value = 1

# Row 10: D10 (Total Attendance Count) and E10 (Real) set to 1
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 1

# Row 11: D11 (Total Attendance Count) and E11 (Real) set to 1
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 1

# Column H (Absent) set to 1 for rows 12-13
$ws.Range("H12:H13").Value = 1

# Row 14: G14 (Invalid) and H14 (Absent) set to 1
$ws.Range("G14").Value = 1
$ws.Range("H14").Value = 1

# Column H (Absent) set to 1 for rows 15-18
$ws.Range("H15:H18").Value = 1
